$d = $word.ActiveDocument

# --- Edit 1 --------------------------------------------------------------
# Table 1, row "data_stb_out" (Output, description column): drop the
# parenthetical "(валидность)" gloss next to "Значимость", leaving just
# "Значимость  данных на шине  ".
$d.Content.Find.Execute(
    "(валидность) данных",   # FindText
    $false,                  # MatchCase
    $false,                  # MatchWholeWord
    $false,                  # MatchWildcards
    $false,                  # MatchSoundsLike
    $false,                  # MatchAllWordForms
    $true,                   # Forward
    1,                       # Wrap (wdFindContinue)
    $false,                  # Format
    " данных",               # ReplaceWith
    2                        # Replace (wdReplaceAll)
) | Out-Null

# --- Edit 2 --------------------------------------------------------------
# Table 1, row "data_ack_in" (Input): the description cell was empty;
# fill it in with the (translated / corrected) explanation of the signal.
$table = $d.Tables.Item(1)
$descCell = $table.Cell(16, 4)
$descRange = $descCell.Range
$descRange.Collapse(1)
$descRange.InsertAfter("Разрешение на прием данных (память ")
$descRange.Collapse(0)
$descRange.InsertAfter("посылает")
$descRange.Collapse(0)
$descRange.InsertAfter(" сигнал ядру что данные актуальны )")
